$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 12 with a copy of row 11's current (pre-edit) data,
# preserving values; column D additionally keeps its custom date format.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols) {
    $src = $ws.Range($col + "11")
    $dst = $ws.Range($col + "12")
    $dst.Value = $src.Value2
}
$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat

# Row 10: newer market date and updated volume
$ws.Range("D10").Value = 44776
$ws.Range("J10").Value = 80

# Row 11: newer market date and updated volume
$ws.Range("D11").Value = 44769
$ws.Range("J11").Value = 50
